# Applies the Mon Sep  9 21:48:29 UTC 2024 crypto-price refresh to the
# "cryptos" worksheet: updates Price (column D) and Volume(1h) (column E)
# for each coin row, and swaps the Polygon/VeChain rows (49/50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value, straight from the diff
$updates = @(
    @{ Cell = 'D2'; Value = '57.316.69' }
    @{ Cell = 'E2'; Value = '  +5.53%  ' }
    @{ Cell = 'D3'; Value = '2.364.37' }
    @{ Cell = 'E3'; Value = '  +4.16%  ' }
    @{ Cell = 'D4'; Value = '1.00' }
    @{ Cell = 'E4'; Value = '  -0.08%  ' }
    @{ Cell = 'D5'; Value = '520.99' }
    @{ Cell = 'E5'; Value = '  +4.34%  ' }
    @{ Cell = 'D6'; Value = '134.78' }
    @{ Cell = 'E6'; Value = '  +4.31%  ' }
    @{ Cell = 'D7'; Value = '0.999' }
    @{ Cell = 'E7'; Value = '  +0.16%  ' }
    @{ Cell = 'D8'; Value = '0.539' }
    @{ Cell = 'E8'; Value = '  +2.46%  ' }
    @{ Cell = 'D9'; Value = '2.359.65' }
    @{ Cell = 'E9'; Value = '  +3.47%  ' }
    @{ Cell = 'E10'; Value = '  +9.00%  ' }
    @{ Cell = 'E11'; Value = '  +0.95%  ' }
    @{ Cell = 'D12'; Value = '5.20' }
    @{ Cell = 'E12'; Value = '  +5.75%  ' }
    @{ Cell = 'D13'; Value = '0.345' }
    @{ Cell = 'E13'; Value = '  +2.52%  ' }
    @{ Cell = 'D14'; Value = '24.02' }
    @{ Cell = 'E14'; Value = '  +3.73%  ' }
    @{ Cell = 'D15'; Value = '2.758.87' }
    @{ Cell = 'E15'; Value = '  +3.22%  ' }
    @{ Cell = 'D16'; Value = '57.112.25' }
    @{ Cell = 'E16'; Value = '  +5.18%  ' }
    @{ Cell = 'E17'; Value = '  +4.83%  ' }
    @{ Cell = 'D18'; Value = '2.351.92' }
    @{ Cell = 'E18'; Value = '  +3.12%  ' }
    @{ Cell = 'D19'; Value = '10.61' }
    @{ Cell = 'E19'; Value = '  +3.19%  ' }
    @{ Cell = 'D20'; Value = '4.29' }
    @{ Cell = 'E20'; Value = '  +3.36%  ' }
    @{ Cell = 'D21'; Value = '322.60' }
    @{ Cell = 'E21'; Value = '  +5.91%  ' }
    @{ Cell = 'D22'; Value = '6.65' }
    @{ Cell = 'E22'; Value = '  +5.21%  ' }
    @{ Cell = 'D23'; Value = '0.997' }
    @{ Cell = 'E23'; Value = '  -0.21%  ' }
    @{ Cell = 'D24'; Value = '61.33' }
    @{ Cell = 'E24'; Value = '  +1.11%  ' }
    @{ Cell = 'D25'; Value = '0.994' }
    @{ Cell = 'E25'; Value = '  -0.26%  ' }
    @{ Cell = 'E26'; Value = '  +7.33%  ' }
    @{ Cell = 'D27'; Value = '7.79' }
    @{ Cell = 'E27'; Value = '  +6.05%  ' }
    @{ Cell = 'D28'; Value = '172.15' }
    @{ Cell = 'E28'; Value = '  -1.78%  ' }
    @{ Cell = 'D29'; Value = '0.0₃0745' }
    @{ Cell = 'E29'; Value = '  +6.20%  ' }
    @{ Cell = 'D30'; Value = '1.21' }
    @{ Cell = 'E30'; Value = '  +11.61%  ' }
    @{ Cell = 'D31'; Value = '6.33' }
    @{ Cell = 'E31'; Value = '  +5.19%  ' }
    @{ Cell = 'D32'; Value = '1.70' }
    @{ Cell = 'E32'; Value = '  +5.33%  ' }
    @{ Cell = 'D33'; Value = '18.44' }
    @{ Cell = 'E33'; Value = '  +3.54%  ' }
    @{ Cell = 'D34'; Value = '0.999' }
    @{ Cell = 'E34'; Value = '  +0.01%  ' }
    @{ Cell = 'D35'; Value = '0.956' }
    @{ Cell = 'E35'; Value = '  +0.38%  ' }
    @{ Cell = 'D36'; Value = '0.999' }
    @{ Cell = 'E36'; Value = '  +0.21%  ' }
    @{ Cell = 'E37'; Value = '  +5.42%  ' }
    @{ Cell = 'D38'; Value = '4.04' }
    @{ Cell = 'E38'; Value = '  +8.33%  ' }
    @{ Cell = 'D39'; Value = '1.52' }
    @{ Cell = 'E39'; Value = '  +8.57%  ' }
    @{ Cell = 'D40'; Value = '37.62' }
    @{ Cell = 'E40'; Value = '  +4.09%  ' }
    @{ Cell = 'E41'; Value = '  +2.06%  ' }
    @{ Cell = 'D42'; Value = '141.02' }
    @{ Cell = 'E42'; Value = '  +12.72%  ' }
    @{ Cell = 'D43'; Value = '3.61' }
    @{ Cell = 'E43'; Value = '  +6.64%  ' }
    @{ Cell = 'D44'; Value = '5.26' }
    @{ Cell = 'E44'; Value = '  +8.89%  ' }
    @{ Cell = 'D45'; Value = '278.62' }
    @{ Cell = 'E45'; Value = '  +13.45%  ' }
    @{ Cell = 'D47'; Value = '0.0932' }
    @{ Cell = 'E47'; Value = '  +3.76%  ' }
    @{ Cell = 'D48'; Value = '0.567' }
    @{ Cell = 'E48'; Value = '  +3.71%  ' }
    @{ Cell = 'B49'; Value = 'Polygon' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Cell = 'D49'; Value = '0.383' }
    @{ Cell = 'E49'; Value = '  +2.13%  ' }
    @{ Cell = 'B50'; Value = 'VeChain' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D50'; Value = '0.0217' }
    @{ Cell = 'E50'; Value = '  +5.78%  ' }
    @{ Cell = 'D51'; Value = '17.06' }
    @{ Cell = 'E51'; Value = '  +4.92%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    # Price strings like "1.00", "0.999" or "5.20" round-trip through the COM
    # Value setter as numbers (dropping the trailing zero / formatting) unless
    # the cell is pre-formatted as Text; force text for anything that looks
    # like a plain number, write it, then drop back to the default style so
    # no stray number format is left behind on the cell.
    if ($u.Value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
